$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value2 = 46077
$ws.Range("C3").Value2 = 46077
$ws.Range("C4").Value2 = 46077
$ws.Range("A5").Value2 = 'A 393-2025'
$ws.Range("B5").Value2 = 45663
$ws.Range("C5").Value2 = 46077
$ws.Range("G5").Value2 = 5.1
$ws.Range("H5").Value2 = 1
$ws.Range("I5").Value2 = 0
$ws.Range("J5").Value2 = 1
$ws.Range("O5").Value2 = 1
$ws.Range("R5").Value2 = 'Spillkråka'
$ws.Range("S5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1481/artfynd/A 393-2025 artfynd.xlsx", "A 393-2025")'
$ws.Range("T5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1481/kartor/A 393-2025 karta.png", "A 393-2025")'
$ws.Range("V5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1481/klagomål/A 393-2025 FSC-klagomål.docx", "A 393-2025")'
$ws.Range("W5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1481/klagomålsmail/A 393-2025 FSC-klagomål mail.docx", "A 393-2025")'
$ws.Range("X5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1481/tillsyn/A 393-2025 tillsynsbegäran.docx", "A 393-2025")'
$ws.Range("Y5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1481/tillsynsmail/A 393-2025 tillsynsbegäran mail.docx", "A 393-2025")'
$ws.Range("Z5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1481/fåglar/A 393-2025 prioriterade fågelarter.docx", "A 393-2025")'
$ws.Range("A6").Value2 = 'A 49789-2023'
$ws.Range("B6").Value2 = 45212
$ws.Range("C6").Value2 = 46077
$ws.Range("G6").Value2 = 3.8
$ws.Range("H6").Value2 = 0
$ws.Range("I6").Value2 = 1
$ws.Range("J6").Value2 = 0
$ws.Range("O6").Value2 = 0
$ws.Range("R6").Value2 = 'Trubbfjädermossa'
$ws.Range("S6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1481/artfynd/A 49789-2023 artfynd.xlsx", "A 49789-2023")'
$ws.Range("T6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1481/kartor/A 49789-2023 karta.png", "A 49789-2023")'
$ws.Range("V6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1481/klagomål/A 49789-2023 FSC-klagomål.docx", "A 49789-2023")'
$ws.Range("W6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1481/klagomålsmail/A 49789-2023 FSC-klagomål mail.docx", "A 49789-2023")'
$ws.Range("X6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1481/tillsyn/A 49789-2023 tillsynsbegäran.docx", "A 49789-2023")'
$ws.Range("Y6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1481/tillsynsmail/A 49789-2023 tillsynsbegäran mail.docx", "A 49789-2023")'
$ws.Range("Z6").ClearContents()
$ws.Range("C7").Value2 = 46077
$ws.Range("C8").Value2 = 46077
$ws.Range("C9").Value2 = 46077
$ws.Range("C10").Value2 = 46077
$ws.Range("C11").Value2 = 46077
$ws.Range("C12").Value2 = 46077
$ws.Range("C13").Value2 = 46077
$ws.Range("C14").Value2 = 46077
$ws.Range("C15").Value2 = 46077
$ws.Range("C16").Value2 = 46077
$ws.Range("A17").Value2 = 'A 8639-2023'
$ws.Range("B17").Value2 = 44977.95614583333
$ws.Range("C17").Value2 = 46077
$ws.Range("G17").Value2 = 3.8
$ws.Range("A18").Value2 = 'A 6258-2024'
$ws.Range("B18").Value2 = 45337.77947916667
$ws.Range("C18").Value2 = 46077
$ws.Range("G18").Value2 = 1.4
$ws.Range("A19").Value2 = 'A 60809-2024'
$ws.Range("B19").Value2 = 45644.61414351852
$ws.Range("C19").Value2 = 46077
$ws.Range("G19").Value2 = 0.5
$ws.Range("A20").Value2 = 'A 11256-2024'
$ws.Range("B20").Value2 = 45371.66233796296
$ws.Range("C20").Value2 = 46077
$ws.Range("G20").Value2 = 0.6
$ws.Range("A21").Value2 = 'A 38-2022'
$ws.Range("B21").Value2 = 44564.37274305556
$ws.Range("C21").Value2 = 46077
$ws.Range("G21").Value2 = 1
$ws.Range("A22").Value2 = 'A 88-2025'
$ws.Range("B22").Value2 = 45659.46386574074
$ws.Range("C22").Value2 = 46077
$ws.Range("G22").Value2 = 1.4
$ws.Range("A23").Value2 = 'A 43067-2024'
$ws.Range("B23").Value2 = 45567.47446759259
$ws.Range("C23").Value2 = 46077
$ws.Range("G23").Value2 = 1.1
$ws.Range("A24").Value2 = 'A 21536-2024'
$ws.Range("B24").Value2 = 45441.59925925926
$ws.Range("C24").Value2 = 46077
$ws.Range("G24").Value2 = 2.8
$ws.Range("A25").Value2 = 'A 12273-2024'
$ws.Range("B25").Value2 = 45378.47817129629
$ws.Range("C25").Value2 = 46077
$ws.Range("G25").Value2 = 0.9
$ws.Range("A26").Value2 = 'A 3811-2024'
$ws.Range("B26").Value2 = 45321.673125
$ws.Range("C26").Value2 = 46077
$ws.Range("A27").Value2 = 'A 56133-2023'
$ws.Range("B27").Value2 = 45240
$ws.Range("C27").Value2 = 46077
$ws.Range("G27").Value2 = 2
$ws.Range("A28").Value2 = 'A 12977-2025'
$ws.Range("B28").Value2 = 45734.45465277778
$ws.Range("C28").Value2 = 46077
$ws.Range("G28").Value2 = 2.1
$ws.Range("A29").Value2 = 'A 61167-2024'
$ws.Range("B29").Value2 = 45645
$ws.Range("C29").Value2 = 46077
$ws.Range("G29").Value2 = 3
$ws.Range("A30").Value2 = 'A 11261-2024'
$ws.Range("B30").Value2 = 45371.67425925926
$ws.Range("C30").Value2 = 46077
$ws.Range("G30").Value2 = 2.7
$ws.Range("A31").Value2 = 'A 19295-2025'
$ws.Range("B31").Value2 = 45769.56212962963
$ws.Range("C31").Value2 = 46077
$ws.Range("G31").Value2 = 2.8
$ws.Range("C32").Value2 = 46077
$ws.Range("A33").Value2 = 'A 24771-2023'
$ws.Range("B33").Value2 = 45084.64277777778
$ws.Range("C33").Value2 = 46077
$ws.Range("G33").Value2 = 1
$ws.Range("A34").Value2 = 'A 44926-2025'
$ws.Range("B34").Value2 = 45918.55856481481
$ws.Range("C34").Value2 = 46077
$ws.Range("G34").Value2 = 3.8
$ws.Range("A35").Value2 = 'A 62804-2023'
$ws.Range("B35").Value2 = 45271
$ws.Range("C35").Value2 = 46077
$ws.Range("G35").Value2 = 0.6
$ws.Range("A36").Value2 = 'A 22953-2023'
$ws.Range("B36").Value2 = 45072
$ws.Range("C36").Value2 = 46077
$ws.Range("G36").Value2 = 1.9
$ws.Range("A37").Value2 = 'A 31321-2025'
$ws.Range("B37").Value2 = 45833
$ws.Range("C37").Value2 = 46077
$ws.Range("G37").Value2 = 6.4
$ws.Range("A38").Value2 = 'A 48265-2025'
$ws.Range("B38").Value2 = 45933
$ws.Range("C38").Value2 = 46077
$ws.Range("G38").Value2 = 2.1
$ws.Range("A39").Value2 = 'A 21477-2025'
$ws.Range("B39").Value2 = 45782.59390046296
$ws.Range("C39").Value2 = 46077
$ws.Range("G39").Value2 = 2.9
$ws.Range("A40").Value2 = 'A 22033-2025'
$ws.Range("B40").Value2 = 45785.29449074074
$ws.Range("C40").Value2 = 46077
$ws.Range("G40").Value2 = 5.2
$ws.Range("A41").Value2 = 'A 62768-2025'
$ws.Range("B41").Value2 = 46008.59856481481
$ws.Range("C41").Value2 = 46077
$ws.Range("G41").Value2 = 4.2
$ws.Range("A42").Value2 = 'A 61178-2024'
$ws.Range("B42").Value2 = 45645
$ws.Range("C42").Value2 = 46077
$ws.Range("G42").Value2 = 6.6
$ws.Range("A43").Value2 = 'A 46379-2025'
$ws.Range("B43").Value2 = 45925
$ws.Range("C43").Value2 = 46077
$ws.Range("G43").Value2 = 7.1
$ws.Range("A44").Value2 = 'A 46384-2025'
$ws.Range("B44").Value2 = 45925
$ws.Range("C44").Value2 = 46077
$ws.Range("G44").Value2 = 3.1
$ws.Range("A45").Value2 = 'A 6042-2024'
$ws.Range("B45").Value2 = 45336
$ws.Range("C45").Value2 = 46077
$ws.Range("G45").Value2 = 1.7
$ws.Range("A46").Value2 = 'A 33201-2023'
$ws.Range("B46").Value2 = 45127.42379629629
$ws.Range("C46").Value2 = 46077
$ws.Range("G46").Value2 = 0.9
$ws.Range("A47").Value2 = 'A 37570-2025'
$ws.Range("B47").Value2 = 45880.37358796296
$ws.Range("C47").Value2 = 46077
$ws.Range("G47").Value2 = 0.9
